# Rename the three header/footer logo pictures (wp:docPr / shape "Name"):
#   - footer Pearson logos:  image2.png -> image1.png   (two occurrences)
#   - header BTEC logo:      image1.jpg -> image2.jpg
#
# These pictures are inline drawings living in the document's header and
# footer stories, so they must be reached through Sections(i).Headers/
# Footers(j).Range.InlineShapes rather than Document.InlineShapes (which
# only covers the main body).

$d = $word.ActiveDocument

function Rename-InlineLogo($shape, $newName) {
    # Setting .Name directly works for shapes that live in the main body
    # or in a header story. Shapes anchored in a footer story raise a
    # stale-handle error when renamed in place, so for those we bounce
    # the picture through a floating Shape and back to an InlineShape,
    # which applies the rename cleanly while preserving the inline
    # (wp:inline) layout.
    try {
        $shape.Name = $newName
    } catch {
        $floating = $shape.ConvertToShape()
        $floating.Name = $newName
        [void]$floating.ConvertToInlineShape()
    }
}

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections($s)

    for ($h = 1; $h -le $section.Headers.Count; $h++) {
        $header = $section.Headers($h)
        if ($header.Exists) {
            for ($i = 1; $i -le $header.Range.InlineShapes.Count; $i++) {
                $shp = $header.Range.InlineShapes($i)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    Rename-InlineLogo $shp "image2.jpg"
                }
            }
        }
    }

    for ($f = 1; $f -le $section.Footers.Count; $f++) {
        $footer = $section.Footers($f)
        if ($footer.Exists) {
            for ($i = 1; $i -le $footer.Range.InlineShapes.Count; $i++) {
                $shp = $footer.Range.InlineShapes($i)
                if ($shp.AlternativeText -like "*PearsonLogo.png") {
                    Rename-InlineLogo $shp "image1.png"
                }
            }
        }
    }
}
